# Automatische test-sync: 2025-06-22 18:48:50
#
# Appends two newly-received mail-log entries to the "Logs" sheet,
# refreshes the "Dashboard" pivot-style summary (re-sorted, with the new
# "Overig" bucket folded in), and extends the bar chart's series ranges to
# cover the Dashboard table's new extent.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append the two new rows (17 & 18)
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A17").Value = "MVO-beleid"
$logs.Range("B17").Value = "mailmind.test@zohomail.eu"
$logs.Range("C17").Value = "Hebben jullie een duurzaamheidsbeleid of MVO-doelen?"
$logs.Range("D17").Value = "Overig"
$logs.Range("F17").Value = "2025-06-22 18:47:57"
$logs.Range("G17").Value = "Nee"

$logs.Range("A18").Value = "Bereikbaarheid klantenservice"
$logs.Range("B18").Value = "mailmind.test@zohomail.eu"
$logs.Range("C18").Value = "Hoe kan ik jullie het beste bereiken?"
$logs.Range("D18").Value = "Overig"
$logs.Range("F18").Value = "2025-06-22 18:48:30"
$logs.Range("G18").Value = "Nee"

# Conditional formatting ranges need to grow from row 16 -> row 18 to keep
# covering the whole Categorie/Beantwoord columns.
$dFormats = $logs.Range("D2:D16").FormatConditions
for ($i = 1; $i -le $dFormats.Count; $i++) {
    $dFormats.Item($i).ModifyAppliesToRange($logs.Range("D2:D18"))
}

$gFormats = $logs.Range("G2:G16").FormatConditions
for ($i = 1; $i -le $gFormats.Count; $i++) {
    $gFormats.Item($i).ModifyAppliesToRange($logs.Range("G2:G18"))
}

# ---------------------------------------------------------------------
# 2. Dashboard sheet: recompute the Categorie/Aantal summary table.
#    The two new "Overig" mails push that category to a count of 2,
#    bumping it above the single-count categories, which shifts the
#    ordering below it down by one row.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A5").Value = "Overig"
$dash.Range("B5").Value = 2

$dash.Range("A6").Value = "Klacht / Probleem"
$dash.Range("B6").Value = 1

$dash.Range("A7").Value = "Uitnodiging / Evenement"
$dash.Range("B7").Value = 1

$dash.Range("A8").Value = "Openingstijden / Locatie"
$dash.Range("B8").Value = 1

$dash.Range("A9").Value = "Samenwerking / Partnerverzoek"
$dash.Range("B9").Value = 1

$dash.Range("A10").Value = "Offerte / Prijsaanvraag"
$dash.Range("B10").Value = 1

$dash.Range("A11").Value = "Retour / Terugbetaling"
$dash.Range("B11").Value = 1

$dash.Range("A12").Value = "Afmelding / Nieuwsbrief"
$dash.Range("B12").Value = 1

# ---------------------------------------------------------------------
# 3. Chart: the bar chart's category/value series referenced
#    Dashboard!$A$2:$A$11 / $B$2:$B$11 - extend by the one extra row.
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$12,'Dashboard'!`$B`$2:`$B`$12,1)"
